$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------------------------------
# Two new metric groups are being added to the "Aged Care" section of the table:
#   - "# Aged Care Active Resident Cases (Weekly)" (+ per-1M / % change variants)
#   - "# Aged Care Active Staff Cases (Weekly)" (+ per-1M / % change variants)
#
# They are slotted into the existing "Aged Care" block (rows 52-61), which bumps the old
# "# Aged Care Active Outbreaks..." rows and the "Treatments" rows that followed them down
# and out of the block; those six displaced rows are re-appended at the bottom of the table
# (new rows 63-68), and the table/autofilter range is grown to match.
# ------------------------------------------------------------------------------------------

# --- New "Resident" weekly-active metrics, replacing the old "# Aged Care Staff Cases" trio ---
$ws.Cells.Item(52,3).Value = "# Aged Care Active Resident Cases (Weekly)"
$ws.Cells.Item(53,3).Value = "# Aged Care Active Resident Cases (Weekly) per 1M"
$ws.Cells.Item(54,3).Value = "% Aged Care Active Resident Cases (Weekly) Change"

# --- "# Aged Care Staff Cases" metrics shift down into rows 55-58 ---
$ws.Cells.Item(55,3).Value = "# Aged Care Staff Cases"
$ws.Cells.Item(56,3).Value = "# Aged Care Staff Cases (7-day avg)"
$ws.Cells.Item(57,3).Value = "# Aged Care Staff Cases (7-day avg) per 1M"
$ws.Cells.Item(58,3).Value = "% Aged Care Staff Cases Weekly Change"

# Rows 56-58 previously used the default cell style for column C; bring them into line with
# the alternate style (fontId 18) used throughout the rest of this "Aged Care" block.
$ws.Cells.Item(48,3).Copy()
$ws.Range("C56:C58").PasteSpecial(-4122)

# --- New "Staff" weekly-active metrics, in rows 59-61 (also adopting the same column C style) ---
$ws.Cells.Item(59,1).Value = "Aged Care"
$ws.Cells.Item(59,2).Value = 60
$ws.Cells.Item(59,3).Value = "# Aged Care Active Staff Cases (Weekly)"

$ws.Cells.Item(60,1).Value = "Aged Care"
$ws.Cells.Item(60,2).Value = 60
$ws.Cells.Item(60,3).Value = "# Aged Care Active Staff Cases (Weekly) per 1M"

$ws.Cells.Item(61,1).Value = "Aged Care"
$ws.Cells.Item(61,2).Value = 60
$ws.Cells.Item(61,3).Value = "% Aged Care Active Staff Cases (Weekly) Change"

$ws.Cells.Item(48,3).Copy()
$ws.Range("C59:C61").PasteSpecial(-4122)

# --- "# Aged Care Active Outbreaks" (first row only) now lands on row 62 ---
$ws.Cells.Item(62,1).Value = "Aged Care"
$ws.Cells.Item(62,2).Value = 60
$ws.Cells.Item(62,3).Value = "# Aged Care Active Outbreaks"

# --- Append 6 new rows (63-68), restoring the metrics displaced off the end of the block ---
$ws.Cells.Item(63,1).Value = "Aged Care"
$ws.Cells.Item(63,2).Value = 60
$ws.Cells.Item(63,3).Value = "# Aged Care Active Outbreaks (7-day avg)"
$ws.Cells.Item(63,4).Value = 560
$ws.Cells.Item(63,6).Value = "X"

$ws.Cells.Item(64,1).Value = "Aged Care"
$ws.Cells.Item(64,2).Value = 60
$ws.Cells.Item(64,3).Value = "# Aged Care Active Outbreaks (7-day avg) per 1M"
$ws.Cells.Item(64,4).Value = 570
$ws.Cells.Item(64,6).Value = "X"

$ws.Cells.Item(65,1).Value = "Aged Care"
$ws.Cells.Item(65,2).Value = 60
$ws.Cells.Item(65,3).Value = "% Aged Care Active Outbreaks Weekly Change"
$ws.Cells.Item(65,4).Value = 580
$ws.Cells.Item(65,6).Value = "X"

$ws.Cells.Item(66,1).Value = "Treatments"
$ws.Cells.Item(66,2).Value = 70
$ws.Cells.Item(66,3).Value = "# Monthly PBS Scripts"
$ws.Cells.Item(66,4).Value = 590
$ws.Cells.Item(66,6).Value = "X"

$ws.Cells.Item(67,1).Value = "Treatments"
$ws.Cells.Item(67,2).Value = 70
$ws.Cells.Item(67,3).Value = "# Monthly PBS Scripts per 1M"
$ws.Cells.Item(67,4).Value = 600
$ws.Cells.Item(67,6).Value = "X"

$ws.Cells.Item(68,1).Value = "Treatments"
$ws.Cells.Item(68,2).Value = 70
$ws.Cells.Item(68,3).Value = "% Monthly PBS Scripts Change"
$ws.Cells.Item(68,4).Value = 610
$ws.Cells.Item(68,6).Value = "X"

# New rows inherit the plain (un-styled) column C look plus the existing E/F formatting used
# by the row directly above them.
$ws.Range("E62:F62").Copy()
$ws.Range("E63:F68").PasteSpecial(-4122)

# --- Expand the table (ListObject) and autofilter to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F68"))

# --- Update selection / scroll position to reflect where the edits were made ---
$ws.Range("C59:C61").Select()
$excel.ActiveWindow.ScrollRow = 34
